# Christian - Updated Chart to handle user input and added information to Introduction page
#
# Adds a "population" column (D) and computes a population-based homeless
# percentage column "pop_percent" (E) next to the existing count-based
# "homeless_percent" column (renamed from "percent", now column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Population figures per county (row -> population) ----
$population = @{
    2  = 19983
    3  = 22582
    4  = 204390
    5  = 77200
    6  = 77331
    7  = 488241
    8  = 3985
    9  = 110593
    10 = 43429
    11 = 7627
    12 = 95222
    13 = 2225
    14 = 97733
    15 = 75061
    16 = 85141
    17 = 32221
    18 = 2252782
    19 = 271473
    20 = 47935
    21 = 22425
    22 = 80707
    23 = 10939
    24 = 66768
    25 = 42243
    26 = 22471
    27 = 13724
    28 = 904980
    29 = 17582
    30 = 129205
    31 = 12083
    32 = 822083
    33 = 522798
    34 = 45723
    35 = 290536
    36 = 4488
    37 = 60760
    38 = 229247
    39 = 50104
    40 = 250873
}

# ---- Header row ----
$ws.Range("C1").Value = "homeless_percent"
$ws.Range("D1").Value = "population"
$ws.Range("E1").Value = "pop_percent"

# ---- Fill in population (D) and pop_percent (E) for every county row ----
foreach ($r in $population.Keys) {
    $ws.Cells.Item($r, 4).Value = $population[$r]
    $ws.Cells.Item($r, 5).Formula = "=B$r/D$r"
    $ws.Cells.Item($r, 5).NumberFormat = "0.00%"
}

# ---- Totals row (41) ----
$ws.Range("D41").Formula = "=SUM(D2:D40)"

# Bold totals, matching the existing bold style already used on row 41.
# (Bold is applied before the number format on each cell so the engine
# doesn't have to create, then abandon, an intermediate un-bolded style.)
$ws.Range("B41").Font.Bold = $true
$ws.Range("B41").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

$ws.Range("C41").Font.Bold = $true
$ws.Range("C41").NumberFormat = "0%"

$ws.Range("D41").Font.Bold = $true
$ws.Range("D41").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# ---- Column width for the new population column ----
$ws.Columns.Item(4).ColumnWidth = 12.17

# ---- View / selection cleanup ----
$ws.Range("F3").Select() | Out-Null

Write-Output "done"
